# Applies the "updated summary charts and summary reports including
# comments from Prof. Erhardt" revision to the Wichita, KS Metro Area-Bus
# FAC summary workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# 1) Workbook window geometry (bookViews/workbookView) - best effort;
#    headless hosts may not persist real screen geometry, but set it
#    via the ActiveWindow object anyway.
# ---------------------------------------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.Height = 15840
    $win.Width  = 29040
    $win.Top    = -120
    $win.Left   = 28680
} catch {
}

# ---------------------------------------------------------------------
# 2) Sheet1 header: Year 1 changed from 2007 to 2012 (two places)
# ---------------------------------------------------------------------
$ws1.Range("C1").Value = 2012
$ws1.Range("E7").Value = 2012

# ---------------------------------------------------------------------
# 3) Sheet1 data table (rows 8-18): updated averages / ridership effect
#    values, reformatted as Number (#,##0.00) for E/F/H and Percent
#    (0.00%) for G/I, and formulas simplified from "*100/X" to "/X"
#    (since the cells are now natively percent-formatted).
# ---------------------------------------------------------------------

# -- Row 8: Vehicle Revenue Miles --
$ws1.Range("E8").Value = 1519833
$ws1.Range("F8").Value = 1749377
$ws1.Range("H8").Value = 137685.47235

# -- Row 9: Average Fare (2018$) --
$ws1.Range("E9").Value = 1.019999059
$ws1.Range("F9").Value = 0.926300675
$ws1.Range("H9").Value = -76487.655283

# -- Row 10: Population + Employment --
$ws1.Range("E10").Value = 914892.67
$ws1.Range("F10").Value = 940329.3199999999
$ws1.Range("H10").Value = 22124.4008985

# -- Row 11: % of Population in Transit Supportive Density --
$ws1.Range("E11").Value = 7.902050484
$ws1.Range("F11").Value = 7.208428975
$ws1.Range("H11").Value = -4725.985050609999

# -- Row 12: Average Gas Price (2018$) --
$ws1.Range("E12").Value = 3.9349
$ws1.Range("F12").Value = 2.72
$ws1.Range("H12").Value = -106844.81968

# -- Row 13: Median Per Capita (2018$) --
$ws1.Range("E13").Value = 28756.38
$ws1.Range("F13").Value = 29196.75
$ws1.Range("H13").Value = 2587.111999999996

# -- Row 14: % of Households with 0 Vehicles --
$ws1.Range("E14").Value = 5.65
$ws1.Range("F14").Value = 6.28
$ws1.Range("H14").Value = 2798.7147904

# -- Row 15: % Working at Home --
$ws1.Range("E15").Value = 3.4
$ws1.Range("F15").Value = 2.875
$ws1.Range("H15").Value = -2769.3364705

# -- Row 16: Years Since Ride-hail Start (E/F stay blank, H unchanged) --
$ws1.Range("H16").Value = -218871.24739

# -- Row 17: Bike Share (E/F unchanged, H unchanged) --
$ws1.Range("E17").Value = 0
$ws1.Range("F17").Value = 1
$ws1.Range("H17").Value = 18208.35942

# -- Row 18: Electric Scooters (E/F unchanged, H unchanged) --
$ws1.Range("E18").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("H18").Value = 0

# -- Row 19: New Reporters - H19 now carries an explicit 0 --
$ws1.Range("H19").Value = 0

# Re-write the "% Diff" formulas for rows 8-21 (percentage sign now
# comes from the cell format, so the "*100" is dropped).
for ($r = 8; $r -le 21; $r++) {
    $ws1.Range("G$r").Formula = "=IFERROR((F$r-E$r)/E$r,0)"
}
# Re-write the "Riddership Effect % Diff" formulas for rows 8-19.
for ($r = 8; $r -le 19; $r++) {
    $ws1.Range("I$r").Formula = "=IFERROR(H$r/`$E`$21,0)"
}

# Apply the "#,##0.00" number format to the value columns.
$ws1.Range("E8:E21").NumberFormat = "#,##0.00"
$ws1.Range("F8:F21").NumberFormat = "#,##0.00"
$ws1.Range("H8:H21").NumberFormat = "#,##0.00"

# Apply the "0.00%" number format to the percent-difference columns.
$ws1.Range("G8:G21").NumberFormat = "0.00%"
$ws1.Range("I8:I21").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 4) Row 20: Total Modeled Ridership
# ---------------------------------------------------------------------
$ws1.Range("E20").Value = 1785748.085
$ws1.Range("F20").Value = 1752757.951
$ws1.Range("G20").Formula = "=IFERROR((F20-E20)/E20,0)"
$ws1.Range("I20").Formula = "=G20"

# ---------------------------------------------------------------------
# 5) Row 21: Total Observed Ridership
# ---------------------------------------------------------------------
$ws1.Range("E21").Value = 1933525
$ws1.Range("F21").Value = 1411363
$ws1.Range("G21").Formula = "=IFERROR((F21-E21)/E21,0)"
$ws1.Range("I21").Formula = "=G21"

# ---------------------------------------------------------------------
# 6) View state: scroll position reset (topLeftCell cleared) and the
#    active selection moved from K20 to H21.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H21").Select()

# ---------------------------------------------------------------------
# 7) Sheet2: "Average Values" label re-centered.
# ---------------------------------------------------------------------
$ws2.Range("E5").HorizontalAlignment = -4108
